# JustinHours.xlsx edit — "hours + added IR schematics + updated fret code header"
#
# - W4 rows 6-11: fill in clock-in/clock-out times (previously blank),
#   which recalculates the Time(hrs) formulas + the Total row.
# - W4 F9/F10 get a new "Write up for MIDI over IR" description, F11 gets
#   a new "Write up for fret slave" description, F8 gets the existing
#   "Investigationg MIDI over IR" description.
# - W3/W4 header cells (B2) get renamed from "Justin Week Two Hours" to
#   "Justin Week Three Hours" / "Justin Week Four Hours" respectively.
# - Active-cell selections are updated per sheet to match where the author
#   ended up after editing.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("W2")
$ws3 = $wb.Worksheets.Item("W3")
$ws4 = $wb.Worksheets.Item("W4")

# --- W4: fill in the blank time entries for rows 6-11 -----------------
$ws4.Range("C6").Value = 0.375
$ws4.Range("D6").Value = 0.4861111111111111

$ws4.Range("C7").Value = 0.53125
$ws4.Range("D7").Value = 0.61111111111111116

$ws4.Range("C8").Value = 0.38541666666666669
$ws4.Range("D8").Value = 0.4861111111111111
$ws4.Range("F8").Value = "Investigationg MIDI over IR"

$ws4.Range("C9").Value = 0.52083333333333337
$ws4.Range("D9").Value = 0.66666666666666663
$ws4.Range("F9").Value = "Write up for MIDI over IR"

$ws4.Range("C10").Value = 0.38541666666666669
$ws4.Range("D10").Value = 0.47222222222222221
$ws4.Range("F10").Value = "Write up for MIDI over IR"

$ws4.Range("C11").Value = 0.54166666666666663
$ws4.Range("D11").Value = 0.625
$ws4.Range("F11").Value = "Write up for fret slave"

# Match the target's number format (0.000) for the E6 / E10 formula cells.
$ws4.Range("E6").NumberFormat = "0.000"
$ws4.Range("E10").NumberFormat = "0.000"

# --- Week headers -------------------------------------------------------
$ws3.Range("B2").Value = "Justin Week Three Hours"
$ws4.Range("B2").Value = "Justin Week Four Hours"

# --- Selections (active cell per sheet, matching author's final state) --
$ws2.Activate()
$ws2.Range("F31").Select()

$ws3.Activate()
$ws3.Range("B2:I2").Select()

$ws4.Activate()
$ws4.Range("K16").Select()
